$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row108 = @(454, 258, 185, 6, 5, 107, 191, 0, 0)
$row109 = @(454, 258, 185, 6, 5, 108, 191, 0, 0)

for ($i = 0; $i -lt $row108.Length; $i++) {
    $ws.Cells.Item(108, $i + 1).Value = $row108[$i]
    $ws.Cells.Item(109, $i + 1).Value = $row109[$i]
}
